$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new year columns P (2019) and Q (2020) to the header row (row 4)
$ws.Cells.Item(4, 16).Value = 2019
$ws.Cells.Item(4, 17).Value = 2020

# Copy style from the preceding year cell (O4) so new cells match formatting
$ws.Cells.Item(4, 15).Copy() | Out-Null
$ws.Range("P4:Q4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Row 5 - renewable energy share values
$ws.Cells.Item(5, 16).Value = 35.67
# Q5 stays empty but should carry the same number style as the rest of the row
$ws.Cells.Item(5, 15).Copy() | Out-Null
$ws.Range("P5:Q5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Row 6 - hydropower electricity production values
$ws.Cells.Item(6, 16).Value = 13859.3
$ws.Cells.Item(6, 17).Value = 13979.1
$ws.Cells.Item(6, 15).Copy() | Out-Null
$ws.Range("P6:Q6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = 0

# Update selection to reflect the author's cursor position after edit
$ws.Range("P9").Select() | Out-Null
